# Update "想去人数" (F column) values across the sheets of
# 北京-漫展信息.xlsx to reflect the latest scrape.
#
# Sheet order in this workbook (see workbook.xml):
#   1 = 展览      (Exhibitions)
#   2 = 演出      (Performances)
#   3 = 本地生活   (Local life)
#   4 = 全部类型   (All types / combined)

$wb = $excel.ActiveWorkbook

$sheetExhibitions = $wb.Worksheets.Item(1)
$sheetPerformances = $wb.Worksheets.Item(2)
$sheetLocalLife = $wb.Worksheets.Item(3)
$sheetAll = $wb.Worksheets.Item(4)

# --- Sheet 1: 展览 ---
$sheetExhibitions.Range("F4").Value = 3432
$sheetExhibitions.Range("F5").Value = 232
$sheetExhibitions.Range("F6").Value = 4944
$sheetExhibitions.Range("F8").Value = 324
$sheetExhibitions.Range("F9").Value = 183
$sheetExhibitions.Range("F10").Value = 653
$sheetExhibitions.Range("F12").Value = 59
$sheetExhibitions.Range("F13").Value = 25
$sheetExhibitions.Range("F14").Value = 679
$sheetExhibitions.Range("F18").Value = 153
$sheetExhibitions.Range("F21").Value = 4811
$sheetExhibitions.Range("F22").Value = 33
$sheetExhibitions.Range("F23").Value = 39
$sheetExhibitions.Range("F25").Value = 5947
$sheetExhibitions.Range("F27").Value = 9
$sheetExhibitions.Range("F28").Value = 3207
$sheetExhibitions.Range("F29").Value = 293
$sheetExhibitions.Range("F30").Value = 689
$sheetExhibitions.Range("F31").Value = 4432
$sheetExhibitions.Range("F34").Value = 137
$sheetExhibitions.Range("F35").Value = 917
$sheetExhibitions.Range("F36").Value = 77
$sheetExhibitions.Range("F37").Value = 18
$sheetExhibitions.Range("F39").Value = 828
$sheetExhibitions.Range("F40").Value = 910

# --- Sheet 2: 演出 ---
$sheetPerformances.Range("F3").Value = 40
$sheetPerformances.Range("F4").Value = 18

# --- Sheet 3: 本地生活 ---
$sheetLocalLife.Range("F3").Value = 1099

# --- Sheet 4: 全部类型 ---
$sheetAll.Range("F4").Value = 1099
$sheetAll.Range("F8").Value = 3432
$sheetAll.Range("F9").Value = 232
$sheetAll.Range("F10").Value = 4944
$sheetAll.Range("F12").Value = 324
$sheetAll.Range("F13").Value = 183
$sheetAll.Range("F14").Value = 653
$sheetAll.Range("F16").Value = 59
$sheetAll.Range("F17").Value = 25
$sheetAll.Range("F18").Value = 679
$sheetAll.Range("F21").Value = 40
$sheetAll.Range("F23").Value = 153
$sheetAll.Range("F26").Value = 4811
$sheetAll.Range("F27").Value = 33
$sheetAll.Range("F28").Value = 39
$sheetAll.Range("F30").Value = 5947
$sheetAll.Range("F32").Value = 9
$sheetAll.Range("F33").Value = 3207
$sheetAll.Range("F34").Value = 293
$sheetAll.Range("F35").Value = 689
$sheetAll.Range("F36").Value = 4432
$sheetAll.Range("F38").Value = 18
$sheetAll.Range("F40").Value = 917
$sheetAll.Range("F41").Value = 77
$sheetAll.Range("F42").Value = 18
$sheetAll.Range("F44").Value = 828
$sheetAll.Range("F45").Value = 910

Write-Output "Updated F-column vote counts across all sheets."
